$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("A3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("A2").Select()
